$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ",\s*"
        $count = $parts.Count
        if ($count -gt 1) {
            $newOrder = @($parts[$count - 1]) + $parts[0..($count - 2)]
            $cell.Value2 = [string]::Join(", ", $newOrder)
        }
    }
}
